$wb = $excel.ActiveWorkbook

# --- Sheet: general ---
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 202.4934176329492
$ws.Range("B4").Value = 0.01800012588500977
$ws.Range("B6").Value = 47.83341763294923
$ws.Range("B9").Value = 47.94
$ws.Range("B10").Value = 106.72

# --- Sheet: x ---
$ws = $wb.Worksheets.Item("x")
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = 4
$ws.Range("B8").Value = 8
$ws.Range("B10").Value = 10
$ws.Range("B12").Value = 12
$ws.Range("B13").Value = 7
$ws.Range("B14").Value = 1

# --- Sheet: U ---
$ws = $wb.Worksheets.Item("U")
$ws.Range("B4").Value = 3
$ws.Range("B7").Value = 3

# --- Sheet: TBar ---
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 39.3596524473688
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 35.00919155153804
$ws.Range("B8").Value = 30
$ws.Range("B9").Value = 35.17863316307213
$ws.Range("B10").Value = 37.32144153802307
$ws.Range("B12").Value = 36.74896288328807
$ws.Range("B13").Value = 44.69631807301934
$ws.Range("B14").Value = 42.53774290968479
$ws.Range("B15").Value = 40.16617764908928

# --- Sheet: Q ---
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 343.05
$ws.Range("C8").Value = 362.28
$ws.Range("C9").Value = 342.825
$ws.Range("C10").Value = 361.69
$ws.Range("C11").Value = 336.1650000000006
$ws.Range("C17").Value = 46.91999999999942
$ws.Range("C18").Value = 36.10499999999942
$ws.Range("C19").Value = 34.91499999999942
$ws.Range("C20").Value = 37.48999999999942
$ws.Range("C21").Value = 39.43499999999941
$ws.Range("C22").Value = 119.55
$ws.Range("C23").Value = 116.1599999999989
$ws.Range("C24").Value = 117.2349999999989
$ws.Range("C25").Value = 121.4449999999989
$ws.Range("C26").Value = 120.25
$ws.Range("C32").Value = 175.1549999999999
$ws.Range("C33").Value = 186.0099999999999
$ws.Range("C34").Value = 167.8249999999999
$ws.Range("C35").Value = 181.8799999999999
$ws.Range("C36").Value = 167.4450000000006
$ws.Range("C37").Value = 238.4400000000016
$ws.Range("C38").Value = 249.0050000000016
$ws.Range("C39").Value = 230.1
$ws.Range("C40").Value = 253.1450000000016
$ws.Range("C41").Value = 239.4750000000016
$ws.Range("C42").Value = 187.4749999999983
$ws.Range("C43").Value = 195.3199999999983
$ws.Range("C44").Value = 177.0549999999983
$ws.Range("C45").Value = 185.2149999999983
$ws.Range("C46").Value = 179.1799999999984
$ws.Range("C49").Value = 152.1250000000007
$ws.Range("C52").Value = 233.1049999999992
$ws.Range("C53").Value = 244.6899999999992
$ws.Range("C54").Value = 229.5549999999992
$ws.Range("C55").Value = 242.5349999999992
$ws.Range("C56").Value = 220.3599999999998
$ws.Range("C57").Value = 343.05
$ws.Range("C58").Value = 362.28
$ws.Range("C59").Value = 342.825
$ws.Range("C60").Value = 361.69
$ws.Range("C61").Value = 336.1650000000006
$ws.Range("C62").Value = 187.4749999999983
$ws.Range("C63").Value = 195.3199999999983
$ws.Range("C64").Value = 177.0549999999983
$ws.Range("C65").Value = 185.2149999999983
$ws.Range("C66").Value = 179.1799999999984
$ws.Range("C67").Value = 238.4400000000016
$ws.Range("C68").Value = 249.0050000000016
$ws.Range("C69").Value = 230.1
$ws.Range("C70").Value = 253.1450000000016
$ws.Range("C71").Value = 239.4750000000016

# --- Sheet: L ---
$ws = $wb.Worksheets.Item("L")
$ws.Range("C3").Value = 12.28
$ws.Range("C5").Value = 11.69
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("C31").Value = 0
# --- Sheet: rho (new data rows) ---
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 1
